# Generate Report for Handoff
# Adds a new handed-off file ("4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md") as row 3
# on the Overview sheet, and the matching per-language detail row 3 on the
# "zh-cn" and "de-de" sheets, mirroring the existing row 2 pattern.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-42-13 20:42:27"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md",
    "",
    "",
    "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-13 20:42:24"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md",
    "",
    "",
    "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf",
    "",
    "",
    "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-13 20:42:27"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md",
    "",
    "",
    "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf",
    "",
    "",
    "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf"
)

Write-Host "Handoff report rows added for 4c2112cc-b53a-4e04-9d2c-d5b1abfd772f"
